$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(1002).RowHeight = 14.25
Write-Host ("UsedRange after: " + $ws.UsedRange.Address())
